# Append 45 more rows of test data (rows 102-146) to the
# master-reg_center_machine_device_h sheet, following the same
# repeating pattern (9-row block of regcntr_id/machine_id combinations)
# as the existing rows, with device_id incrementing sequentially.
# Then re-select the rows below the new data and set the page to
# portrait orientation, matching the author's follow-up save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 102
$endRow   = 146
$startC   = 3000121

for ($r = $startRow; $r -le $endRow; $r++) {
    $idx = $r - $startRow
    $a = 10002 + ($idx % 9)
    $b = 10021 + ($idx % 9)
    $c = $startC + $idx

    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin()"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

# Select everything below the newly pasted data (mirrors the
# "select remaining rows" action visible in the saved selection).
[void]$ws.Range("A147:XFD1048576").Select()

# Page setup tweak captured in the same save.
$ws.PageSetup.Orientation = 1
